$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report regeneration:
#   File 1 GUID: 24222ef3-212d-4d3e-83a6-75c7416c571f -> 22b62f1f-4ee9-4abb-b413-b037cdcf3750
#   File 2 GUID: 8c9297c7-c14f-4444-b02e-345cf21e377c -> ffff7717bcbd-2836-48ac-8573-e4baf192c19c
#   Handoff/handback xliff hash (now shared by both files):
#       a7d4875fc4555919c4d354939c6f6863e3b47f77 / b3e3e337c85a8ed4ffac93a40cecc2367efb44fc
#       -> 712d9bdd128fdce9de48f75a30dae591a7aa73d2
#   Plus refreshed timestamps.
# ---------------------------------------------------------------------------

$oldGuid1 = "24222ef3-212d-4d3e-83a6-75c7416c571f"
$newGuid1 = "22b62f1f-4ee9-4abb-b413-b037cdcf3750"
$oldGuid2 = "8c9297c7-c14f-4444-b02e-345cf21e377c"
$newGuid2 = "ffff7717bcbd-2836-48ac-8573-e4baf192c19c"

$newHashZh = "$newGuid1.712d9bdd128fdce9de48f75a30dae591a7aa73d2.zh-cn.xlf"
$newHashDe = "$newGuid1.712d9bdd128fdce9de48f75a30dae591a7aa73d2.de-de.xlf"

function Set-CellText($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

function Update-Hyperlink($ws, $oldText, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -eq $oldText) {
            $hl.TextToDisplay = $newText
        }
    }
}

# ------------------------- Overview sheet ---------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellText $wsOverview "A2" "$newGuid1.md"
Set-CellText $wsOverview "B2" "e2e\$newGuid1.md"
Set-CellText $wsOverview "G2" "2016-08-18 07:02:17"

Set-CellText $wsOverview "A3" "$newGuid2.md"
Set-CellText $wsOverview "B3" "e2e\$newGuid2.md"
Set-CellText $wsOverview "G3" "2016-08-18 07:02:17"

Update-Hyperlink $wsOverview "e2e\$oldGuid1.md" "e2e\$newGuid1.md"
Update-Hyperlink $wsOverview "e2e\$oldGuid2.md" "e2e\$newGuid2.md"

# ------------------------- zh-cn sheet -------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellText $wsZhCn "A2" "$newGuid1.md"
Set-CellText $wsZhCn "G2" $newHashZh
Set-CellText $wsZhCn "H2" "2016-08-18 07:02:11"
Set-CellText $wsZhCn "I2" "$newGuid1.md"
Set-CellText $wsZhCn "J2" $newHashZh
Set-CellText $wsZhCn "K2" "2016-08-18 07:02:38"

Set-CellText $wsZhCn "A3" "$newGuid2.md"
Set-CellText $wsZhCn "G3" $newHashZh
Set-CellText $wsZhCn "H3" "2016-08-18 07:02:11"
Set-CellText $wsZhCn "I3" "$newGuid2.md"
Set-CellText $wsZhCn "J3" $newHashZh
Set-CellText $wsZhCn "K3" "2016-08-18 07:02:38"

Update-Hyperlink $wsZhCn "$oldGuid1.md" "$newGuid1.md"
Update-Hyperlink $wsZhCn "$oldGuid2.md" "$newGuid2.md"

# ------------------------- de-de sheet -------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellText $wsDeDe "A2" "$newGuid1.md"
Set-CellText $wsDeDe "G2" $newHashDe
Set-CellText $wsDeDe "H2" "2016-08-18 07:02:17"
Set-CellText $wsDeDe "I2" "$newGuid1.md"
Set-CellText $wsDeDe "J2" $newHashDe
Set-CellText $wsDeDe "K2" "2016-08-18 07:02:46"

Set-CellText $wsDeDe "A3" "$newGuid2.md"
Set-CellText $wsDeDe "G3" $newHashDe
Set-CellText $wsDeDe "H3" "2016-08-18 07:02:17"
Set-CellText $wsDeDe "I3" "$newGuid2.md"
Set-CellText $wsDeDe "J3" $newHashDe
Set-CellText $wsDeDe "K3" "2016-08-18 07:02:46"

Update-Hyperlink $wsDeDe "$oldGuid1.md" "$newGuid1.md"
Update-Hyperlink $wsDeDe "$oldGuid2.md" "$newGuid2.md"
